$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.646.01"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.565.01"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.83"
$ws.Range("E8").Value = "  +5.17%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.789.19"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "1.564.81"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "28.680.67"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.49"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.88"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.54"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0458"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.405.84"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.03"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.517"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.767"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.90"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").Value = "1.701.66"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.858"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.88"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.70"
$ws.Range("E50").Value = "  +5.41%  "
$ws.Range("E51").Value = "  -0.46%  "
